# Updated symbol list on Fri Dec 16 23:28:45 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-typed numeric-looking value into column D so Excel
# keeps it as a string (matches the original inlineStr cells) instead of
# auto-converting it to a number.
function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "231.38"
Set-TextValue "D3"  "22.77"
Set-TextValue "D4"  "5.560"
Set-TextValue "D5"  "0.05577"
Set-TextValue "D6"  "3.421"
Set-TextValue "D7"  "6.494"
Set-TextValue "D8"  "1.136"
Set-TextValue "D9"  "0.7951"
Set-TextValue "D11" "0.07384"
Set-TextValue "D12" "0.03157"
Set-TextValue "D13" "0.02944"
Set-TextValue "D14" "0.09244"
Set-TextValue "D15" "0.001678"
Set-TextValue "D16" "3.279"
Set-TextValue "D17" "0.04733"
Set-TextValue "D18" "0.0005947"
Set-TextValue "D19" "0.006254"
Set-TextValue "D20" "0.005272"
Set-TextValue "D21" "0.001069"
Set-TextValue "D22" "0.0001509"
Set-TextValue "D23" "3.678"
Set-TextValue "D26" "0.1282"
Set-TextValue "D27" "0.0008359"
Set-TextValue "D40" "0.04056"
Set-TextValue "D41" "0.007184"

# --- Row 18 (One): Volume(1h) label change ---
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 42/43: CEJI and BKEXToken swap places ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1032"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003350"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Row 44 (LocalTraders) ---
Set-TextValue "D44" "0.008152"
$ws.Range("E44").Value = "43LocalTradersLCT"

# --- Remaining price (column D) updates ---
Set-TextValue "D46" "0.00005560"
Set-TextValue "D47" "0.00000000755"
Set-TextValue "D48" "0.6792"
Set-TextValue "D49" "0.09517"
Set-TextValue "D50" "0.00002113"
Set-TextValue "D51" "0.01016"

Write-Output "edits applied"
